$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Task Summary")
$ws2 = $wb.Worksheets.Item("Daily Tracker")

# --- Task Summary: add Task 29 / Task 30 summary rows ---
$ws1.Range("A41").Value = "Task 29"
$ws1.Range("B41").Value = 44456
$ws1.Range("C41").Value = "PGP and SSH keys"
$ws1.Range("D41").Value = 44460

$ws1.Range("A42").Value = "Task 30"
$ws1.Range("B42").Value = 44456
$ws1.Range("C42").Value = "Bank Payment Technical Documentation"
$ws1.Range("D42").Value = 44466
$ws1.Range("E42").Value = "Y"

$ws1.Rows.Item(41).RowHeight = 15.75
$ws1.Rows.Item(42).RowHeight = 15.75

# --- Daily Tracker: update row 69 (was generic "Nil" placeholder) and append new rows ---
$ws2.Range("B69").Value = "Task 29"
$ws2.Range("C69").Value = "PGP Keys"
$ws2.Range("D69").Value = "PGP Keys"
$ws2.Range("E69").Value = "PGP keys analysed"

$ws2.Range("A70").Value = 44460
$ws2.Range("B70").Value = "Task 29"
$ws2.Range("C70").Value = "SSH Keys"
$ws2.Range("D70").Value = "SSH Keys"
$ws2.Range("E70").Value = "SSH keys analysed"

$ws2.Range("A71").Value = 44461
$ws2.Range("B71").Value = "Task 30"
$ws2.Range("C71").Value = "Bank Payment Technical Documentation"
$ws2.Range("D71").Value = "Bank Payment Technical Documentation"
$ws2.Range("E71").Value = "Partially completed"

$ws2.Range("A72").Value = 44462
$ws2.Range("B72").Value = "Task 30"
$ws2.Range("C72").Value = "Bank Payment Technical Documentation"
$ws2.Range("D72").Value = "Bank Payment Technical Documentation"
$ws2.Range("E72").Value = "Partially completed"

$ws2.Range("A73").Value = 44463
$ws2.Range("B73").Value = "Task 30"
$ws2.Range("C73").Value = "Bank Payment Technical Documentation"
$ws2.Range("D73").Value = "Bank Payment Technical Documentation"
$ws2.Range("E73").Value = "Partially completed"

$ws2.Range("A74").Value = 44466
$ws2.Range("B74").Value = "Task 30"
$ws2.Range("C74").Value = "Bank Payment Technical Documentation"
$ws2.Range("D74").Value = "Bank Payment Technical Documentation"
$ws2.Range("E74").Value = "Documentation completed"

$ws2.Rows.Item(70).RowHeight = 15.75
$ws2.Rows.Item(71).RowHeight = 15.75
$ws2.Rows.Item(72).RowHeight = 15.75
$ws2.Rows.Item(73).RowHeight = 15.75
$ws2.Rows.Item(74).RowHeight = 15.75

# --- View state: active sheet switches back to Task Summary, selections updated ---
$null = $ws2.Range("E75").Select()
$null = $ws1.Activate()
$null = $ws1.Range("G66").Select()
